$d = $word.ActiveDocument

# The "KEY ACHIEVEMENTS AND IMPACT" -> "Impact" section originally contains six
# bullet paragraphs (paragraphs 55-60, 1-based). Rewrite them as four
# impact-focused accomplishment statements per the commit.

# 1) Replace the text of the first three bullets in place.
$d.Paragraphs.Item(55).Range.Text = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
$d.Paragraphs.Item(56).Range.Text = "• `$4.7M savings enabled nonprofit access"
$d.Paragraphs.Item(57).Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# 2) Replace the text of the sixth (last) bullet, which survives as the fourth bullet.
$d.Paragraphs.Item(60).Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# 3) Remove the fourth and fifth bullets entirely (their whole paragraphs, including marks).
$d.Paragraphs.Item(59).Range.Delete()
$d.Paragraphs.Item(58).Range.Delete()
